$wb = $excel.ActiveWorkbook

# Sheet AmazGoodr (currently active / tabSelected) -> add a new row 4, copy of row 2 with updated threshold
$wsGoodr = $wb.Worksheets.Item("AmazGoodr")
$wsGoodr.Range("A2:AA2").Copy() | Out-Null
$wsGoodr.Range("A4").PasteSpecial() | Out-Null
$wsGoodr.Range("A4").Value = "threshold 0,7"
$wsGoodr.Range("S4").Value = 0.7
$wsGoodr.Range("W4").Value = "0:04:26.095"
$wsGoodr.Range("X4").Value = 0.9048
$wsGoodr.Range("Y4").Value = 0.5938
$wsGoodr.Range("Z4").Value = 0.717

# Sheet GoodrCov -> add a new row 3, copy of row 2 with updated threshold
$wsCov = $wb.Worksheets.Item("GoodrCov")
$wsCov.Range("A2:AE2").Copy() | Out-Null
$wsCov.Range("A3").PasteSpecial() | Out-Null
$wsCov.Range("A3").Value = "threshold 0,7"
$wsCov.Range("W3").Value = 0.7
$wsCov.Range("AA3").Value = "0:12:20.514"
$wsCov.Range("AB3").Value = 0.9545
$wsCov.Range("AC3").Value = 0.7
$wsCov.Range("AD3").Value = 0.8077

# Make GoodrCov the active/selected sheet, with the indicated selection
$wsCov.Activate()
$wsCov.Application.ActiveWindow.ScrollColumn = 8
$wsCov.Range("AE3").Select() | Out-Null

$wsGoodr.Range("A2").Select() | Out-Null
